$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.250631
$ws.Range("H2").Value = 0.751893
$ws.Range("I2").Value = 0.2648339568266264
$ws.Range("J2").Value = 0.2648339568266264
$ws.Range("M2").Value = 1.672411
$ws.Range("N2").Value = 5.017233
$ws.Range("O2").Value = 0.6245395681653219
$ws.Range("P2").Value = 0.624539568165322
$ws.Range("Q2").Value = 0.419158041341
$ws.Range("R2").Value = 3.772422372069
$ws.Range("S2").Value = 0.1653992850320147
$ws.Range("T2").Value = 0.1653992850320148
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.250631
$ws.Range("H3").Value = 0.751893
$ws.Range("I3").Value = 0.2648339568266264
$ws.Range("J3").Value = 0.2648339568266264
$ws.Range("O3").Value = 0.2513435317223857
$ws.Range("P3").Value = 0.2513435317223857
$ws.Range("Q3").Value = 0.1686885312486667
$ws.Range("R3").Value = 1.518196781238
$ws.Range("S3").Value = 0.06656430202881809
$ws.Range("T3").Value = 0.06656430202881809
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.250631
$ws.Range("H4").Value = 0.751893
$ws.Range("I4").Value = 0.2648339568266264
$ws.Range("J4").Value = 0.2648339568266264
$ws.Range("N4").Value = 0.997092
$ws.Range("O4").Value = 0.1241169001122924
$ws.Range("P4").Value = 0.1241169001122924
$ws.Range("Q4").Value = 0.08330072168399999
$ws.Range("R4").Value = 0.749706495156
$ws.Range("S4").Value = 0.03287036976579354
$ws.Range("T4").Value = 0.03287036976579354
# Row 5
$ws.Range("I5").Value = 0.2480790641859371
$ws.Range("J5").Value = 0.2480790641859371
$ws.Range("M5").Value = 1.672411
$ws.Range("N5").Value = 5.017233
$ws.Range("O5").Value = 0.6245395681653219
$ws.Range("P5").Value = 0.624539568165322
$ws.Range("Q5").Value = 0.3926397350546667
$ws.Range("R5").Value = 3.533757615492001
$ws.Range("S5").Value = 0.1549351916175423
$ws.Range("T5").Value = 0.1549351916175424
# Row 6
$ws.Range("I6").Value = 0.2480790641859371
$ws.Range("J6").Value = 0.2480790641859371
$ws.Range("O6").Value = 0.2513435317223857
$ws.Range("P6").Value = 0.2513435317223857
$ws.Range("S6").Value = 0.06235306813887784
$ws.Range("T6").Value = 0.06235306813887784
# Row 7
$ws.Range("I7").Value = 0.2480790641859371
$ws.Range("J7").Value = 0.2480790641859371
$ws.Range("N7").Value = 0.997092
$ws.Range("O7").Value = 0.1241169001122924
$ws.Range("P7").Value = 0.1241169001122924
$ws.Range("Q7").Value = 0.078030647312
$ws.Range("S7").Value = 0.03079080442951693
$ws.Range("T7").Value = 0.03079080442951693
# Row 8
$ws.Range("G8").Value = 0.4609646666666666
$ws.Range("I8").Value = 0.4870869789874365
$ws.Range("J8").Value = 0.4870869789874365
$ws.Range("M8").Value = 1.672411
$ws.Range("N8").Value = 5.017233
$ws.Range("O8").Value = 0.6245395681653219
$ws.Range("P8").Value = 0.624539568165322
$ws.Range("Q8").Value = 0.7709223791446667
$ws.Range("R8").Value = 6.938301412302
$ws.Range("S8").Value = 0.3042050915157648
$ws.Range("T8").Value = 0.3042050915157649
# Row 9
$ws.Range("G9").Value = 0.4609646666666666
$ws.Range("I9").Value = 0.4870869789874365
$ws.Range("J9").Value = 0.4870869789874365
$ws.Range("O9").Value = 0.2513435317223857
$ws.Range("P9").Value = 0.2513435317223857
$ws.Range("Q9").Value = 0.3102547273782222
$ws.Range("S9").Value = 0.1224261615546898
$ws.Range("T9").Value = 0.1224261615546898
# Row 10
$ws.Range("G10").Value = 0.4609646666666666
$ws.Range("I10").Value = 0.4870869789874365
$ws.Range("J10").Value = 0.4870869789874365
$ws.Range("N10").Value = 0.997092
$ws.Range("O10").Value = 0.1241169001122924
$ws.Range("P10").Value = 0.1241169001122924
$ws.Range("S10").Value = 0.06045572591698192
$ws.Range("T10").Value = 0.06045572591698193
